$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C8 value from 0 to 222 (Integer min for rule R10)
$ws.Range("C8").Value = 222
